$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (pushes old D..K to F..M)
$ws.Columns("D:E").Insert()

# Populate the two new columns (D and E) with updated quarterly figures
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 986900
$ws.Range("E8").Value = 952700
$ws.Range("D9").Value = 653900
$ws.Range("E9").Value = 636000
$ws.Range("D10").Value = 333000
$ws.Range("E10").Value = 316700
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 10600
$ws.Range("E14").Value = 18900
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 893700
$ws.Range("E17").Value = 890500
$ws.Range("D18").Value = 93200
$ws.Range("E18").Value = 62200
$ws.Range("D20").Value = -200
$ws.Range("E20").Value = -4000
$ws.Range("D21").Value = 120300
$ws.Range("E21").Value = 85800
$ws.Range("D22").Value = 14500
$ws.Range("E22").Value = 13800
$ws.Range("D23").Value = 78600
$ws.Range("E23").Value = 44400
$ws.Range("D24").Value = 19900
$ws.Range("E24").Value = 14900
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 58700
$ws.Range("E26").Value = 29400
$ws.Range("D27").Value = 57400
$ws.Range("E27").Value = 28200
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 5700
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 200
$ws.Range("E32").Value = 4000
$ws.Range("D33").Value = 63100
$ws.Range("E33").Value = 28200
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 63100
$ws.Range("E35").Value = 28200
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 619700
$ws.Range("E41").Value = 529900
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 1021000
$ws.Range("E43").Value = 1041800
$ws.Range("D44").Value = 633900
$ws.Range("E44").Value = 655700
$ws.Range("D45").Value = 108600
$ws.Range("E45").Value = 97200
$ws.Range("D46").Value = 2383100
$ws.Range("E46").Value = 2324700
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 610100
$ws.Range("E48").Value = 608700
$ws.Range("D49").Value = 1388200
$ws.Range("E49").Value = 1399600
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 234800
$ws.Range("E52").Value = 272000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 4616300
$ws.Range("E54").Value = 4605100
$ws.Range("D57").Value = 418900
$ws.Range("E57").Value = 400300
$ws.Range("D58").Value = 68200
$ws.Range("E58").Value = 67300
$ws.Range("D59").Value = 593900
$ws.Range("E59").Value = 572500
$ws.Range("D60").Value = 1081000
$ws.Range("E60").Value = 1040100
$ws.Range("D61").Value = 1414800
$ws.Range("E61").Value = 1436700
$ws.Range("D62").Value = 459700
$ws.Range("E62").Value = 497500
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 2974000
$ws.Range("E66").Value = 2992000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 3543000
$ws.Range("E72").Value = 3505100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1642300
$ws.Range("E76").Value = 1613000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 63100
$ws.Range("E81").Value = 28200
$ws.Range("D83").Value = 27300
$ws.Range("E83").Value = 27600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 164500
$ws.Range("E89").Value = 83100
$ws.Range("D91").Value = -34000
$ws.Range("E91").Value = -18200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -31900
$ws.Range("E94").Value = -18700
$ws.Range("D96").Value = -24900
$ws.Range("E96").Value = -24900
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -40100
$ws.Range("E100").Value = -46000
$ws.Range("D101").Value = -2800
$ws.Range("E101").Value = -5900
$ws.Range("D102").Value = 89700
$ws.Range("E102").Value = 12500
